$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Basic Game rubric")
$ws.Activate()

$ws.Range("B5").Value = 3
$ws.Range("B5").Select()
